$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2
$ws.Range("C2").Value = 0.548936170212766
$ws.Range("P2").Value = 0.1617021276595745
$ws.Range("S2").Value = 0.08936170212765958
$ws.Range("C3").Value = 0.03597122302158273
$ws.Range("J3").Value = 0.02877697841726619
$ws.Range("P3").Value = 0.7482014388489209
$ws.Range("S3").Value = 0.1870503597122302
$ws.Range("J4").Value = 0.03846153846153846
$ws.Range("P4").Value = 0.6538461538461539
$ws.Range("S4").Value = 0.3076923076923077
$ws.Range("B6").Value = 0.04017857142857143
$ws.Range("D6").Value = 0.02232142857142857
$ws.Range("F6").Value = 0.07589285714285714
$ws.Range("J6").Value = 0.2455357142857143
$ws.Range("O6").Value = 0.01785714285714286
$ws.Range("Q6").Value = 0.1607142857142857
$ws.Range("R6").Value = 0.07589285714285714
$ws.Range("S6").Value = 0.3616071428571428
$ws.Range("B7").Value = 0.08021390374331551
$ws.Range("D7").Value = 0.0267379679144385
$ws.Range("F7").Value = 0.0748663101604278
$ws.Range("J7").Value = 0.106951871657754
$ws.Range("O7").Value = 0.0213903743315508
$ws.Range("Q7").Value = 0.160427807486631
$ws.Range("R7").Value = 0.1016042780748663
$ws.Range("S7").Value = 0.427807486631016
$ws.Range("B8").Value = 0.08158995815899582
$ws.Range("D8").Value = 0.02301255230125523
$ws.Range("E8").Value = 0.002092050209205021
$ws.Range("F8").Value = 0.04811715481171548
$ws.Range("J8").Value = 0.1066945606694561
$ws.Range("O8").Value = 0.0104602510460251
$ws.Range("Q8").Value = 0.1903765690376569
$ws.Range("R8").Value = 0.09832635983263599
$ws.Range("S8").Value = 0.4393305439330544
$ws.Range("B9").Value = 0.07253886010362694
$ws.Range("D9").Value = 0.0155440414507772
$ws.Range("F9").Value = 0.03626943005181347
$ws.Range("J9").Value = 0.09844559585492228
$ws.Range("O9").Value = 0.01036269430051814
$ws.Range("Q9").Value = 0.1761658031088083
$ws.Range("R9").Value = 0.1036269430051813
$ws.Range("S9").Value = 0.4870466321243523
$ws.Range("B10").Value = 0.08602999210734018
$ws.Range("D10").Value = 0.02525651144435675
$ws.Range("E10").Value = 0.0007892659826361484
$ws.Range("F10").Value = 0.06471981057616416
$ws.Range("J10").Value = 0.1231254932912391
$ws.Range("O10").Value = 0.01262825572217837
$ws.Range("Q10").Value = 0.2067876874506709
$ws.Range("R10").Value = 0.09865824782951854
$ws.Range("S10").Value = 0.3820047355958958
$ws.Range("G11").Value = 0.1559322033898305
$ws.Range("J11").Value = 0.0847457627118644
$ws.Range("K11").Value = 0.2135593220338983
$ws.Range("L11").Value = 0.5322033898305085
$ws.Range("S11").Value = 0.0135593220338983
$ws.Range("G12").Value = 0.7267080745341615
$ws.Range("J12").Value = 0.2049689440993789
$ws.Range("K12").Value = 0.01863354037267081
$ws.Range("L12").Value = 0.0124223602484472
$ws.Range("S12").Value = 0.03726708074534162
$ws.Range("F13").Value = 0.02222222222222222
$ws.Range("G13").Value = 0.7555555555555555
$ws.Range("J13").Value = 0.2222222222222222
$ws.Range("F15").Value = 0.02369668246445497
$ws.Range("H15").Value = 0.2417061611374408
$ws.Range("I15").Value = 0.05687203791469194
$ws.Range("J15").Value = 0.3270142180094787
$ws.Range("K15").Value = 0.07109004739336493
$ws.Range("O15").Value = 0.09004739336492891
$ws.Range("S15").Value = 0.1895734597156398
$ws.Range("F16").Value = 0.01764705882352941
$ws.Range("H16").Value = 0.1705882352941177
$ws.Range("I16").Value = 0.07058823529411765
$ws.Range("J16").Value = 0.4352941176470588
$ws.Range("K16").Value = 0.08235294117647059
$ws.Range("M16").Value = 0.02352941176470588
$ws.Range("O16").Value = 0.07058823529411765
$ws.Range("S16").Value = 0.1294117647058824
$ws.Range("F17").Value = 0.02237136465324385
$ws.Range("H17").Value = 0.203579418344519
$ws.Range("I17").Value = 0.1029082774049217
$ws.Range("J17").Value = 0.4093959731543624
$ws.Range("K17").Value = 0.09619686800894854
$ws.Range("M17").Value = 0.01789709172259508
$ws.Range("N17").Value = 0.002237136465324385
$ws.Range("O17").Value = 0.05369127516778523
$ws.Range("S17").Value = 0.09172259507829977
$ws.Range("F18").Value = 0.02192982456140351
$ws.Range("H18").Value = 0.1710526315789474
$ws.Range("I18").Value = 0.1052631578947368
$ws.Range("J18").Value = 0.3815789473684211
$ws.Range("K18").Value = 0.1052631578947368
$ws.Range("M18").Value = 0.01754385964912281
$ws.Range("O18").Value = 0.07456140350877193
$ws.Range("S18").Value = 0.1228070175438596
$ws.Range("F19").Value = 0.02249806051202483
$ws.Range("H19").Value = 0.2102404965089217
$ws.Range("I19").Value = 0.07525213343677269
$ws.Range("J19").Value = 0.3801396431342126
$ws.Range("K19").Value = 0.1024049650892164
$ws.Range("M19").Value = 0.02249806051202483
$ws.Range("O19").Value = 0.06671838634600466
$ws.Range("S19").Value = 0.1202482544608223
